$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Invalid(G), Absent(H)
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Rows 4-7: Total Attendance Count(D), Real(E)
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1

# Rows 8-9: Absent(H)
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1

# Rows 10-13: Total Attendance Count(D), Real(E)
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1

# Row 14: Absent(H)
$ws.Range("H14").Value = 1

# Row 15: Total Attendance Count(D), Real(E)
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1

# Row 16: Absent(H)
$ws.Range("H16").Value = 1

# Row 17: Total Attendance Count(D), Real(E)
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 1

# Row 18: Absent(H)
$ws.Range("H18").Value = 1
